# fix: save unselect status to db
#
# The reconciliation table gained two new (unreconciled) ledger lines —
# id 1 and id 2, both dated 2024-10-04 — which pushed the previously
# existing rows (ids 4, 8, 13, 15) down by one row. The sheet grows from
# 5 data rows (A1:K6) to 6 data rows (A1:K7).
#
# Rather than using Rows.Insert() (which stamps an inherited style index
# onto the shifted cells), every row 2-7 is rewritten in place with its
# final id/data/historico/lote/lancamento/d-or-c/dc/conta/conciliada
# values, matching the row-by-row shift visible in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E, F, G, H and J hold numeric-looking text (lote, lancamento,
# débito/crédito amount, conta) that must stay shared-string text (as in
# the source file) rather than being coerced to numbers — and must not
# pick up a leading zero loss (conta "02089"). A leading apostrophe forces
# text entry; re-applying the "Normal" style afterwards drops the
# quote-prefix formatting so the cell keeps the workbook's default style.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# id, data, historico, contra_partida, lote, lancamento, d, c, dc, conta, conciliada
$rows = @(
    @{ Row=2; Id=1;  Data="2024-10-04 00:00:00"; Historico="VLR. NF. N. 570  DEUCHER RESTAURANTE EIRELI";     Lote="4572"; Lancamento="10";  D=$null;    C="16550"; Dc="C"; Conta="02089"; Conciliada=0 }
    @{ Row=3; Id=2;  Data="2024-10-04 00:00:00"; Historico="VLR. NF. N. 569  DEUCHER RESTAURANTE EIRELI";     Lote="4572"; Lancamento="12";  D=$null;    C="25165"; Dc="C"; Conta="02089"; Conciliada=0 }
    @{ Row=4; Id=4;  Data="2024-10-08 00:00:00"; Historico="N/PGTO. NF. N. 570  DEUCHER RESTAURANTE EIRELI";  Lote="4828"; Lancamento="103"; D="16550";   C=$null;   Dc="D"; Conta="02089"; Conciliada=0 }
    @{ Row=5; Id=8;  Data="2024-10-23 00:00:00"; Historico="VLR. NF. N. 578  DEUCHER RESTAURANTE EIRELI";     Lote="4974"; Lancamento="12";  D=$null;    C="81640"; Dc="C"; Conta="02089"; Conciliada=0 }
    @{ Row=6; Id=13; Data="2024-10-25 00:00:00"; Historico="N/PGTO. NF. N. 578  DEUCHER RESTAURANTE EIRELI";  Lote="5063"; Lancamento="51";  D="81640";   C=$null;   Dc="C"; Conta="02089"; Conciliada=0 }
    @{ Row=7; Id=15; Data="2024-10-25 00:00:00"; Historico="N/PGTO. NF. N. 580  DEUCHER RESTAURANTE EIRELI";  Lote="5063"; Lancamento="55";  D="50220";   C=$null;   Dc="C"; Conta="02089"; Conciliada=0 }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value = $r.Id
    $ws.Range("B$n").Value = $r.Data
    $ws.Range("C$n").Value = $r.Historico
    $ws.Range("D$n").Value = $null

    Set-TextValue $ws.Range("E$n") $r.Lote
    Set-TextValue $ws.Range("F$n") $r.Lancamento

    if ($r.D -ne $null) { Set-TextValue $ws.Range("G$n") $r.D } else { $ws.Range("G$n").Value = $null }
    if ($r.C -ne $null) { Set-TextValue $ws.Range("H$n") $r.C } else { $ws.Range("H$n").Value = $null }

    $ws.Range("I$n").Value = $r.Dc
    Set-TextValue $ws.Range("J$n") $r.Conta
    $ws.Range("K$n").Value = $r.Conciliada
}

Write-Output "applied $($rows.Count) reconciliation rows (A1:K7)"
